# edit.ps1 - applies the relatorioParecer.docx changes described by the diff:
#   1. "I - DESCRIÇÃO"                   -> "I - " + "tipodocumento"
#   2. "Vem ao exame desta Comissão o "  -> "Vem ao exame desta Comissão o" + " " +
#                                            "pedidoAprovacao" + " "
#   3. ", sem alocação de Carga Horária." -> ", " + "cargaHoraria" + "."
#
# Word merges two adjacent runs back together when it serializes the
# document if they end up with identical run formatting, so a plain
# Range.InsertAfter() alone is not enough to obtain separate <w:r>
# elements. Temporarily anchoring a bookmark exactly at the split point
# forces Word to keep the run boundary; the bookmark must still be present
# at the moment the new text is inserted, and can then be deleted right
# away (it never ends up in the saved document).

$d = $word.ActiveDocument

function Insert-AsNewRun($rng, [string]$text) {
    $bk = "zzzsplit" + [System.Guid]::NewGuid().ToString("N")
    $d.Bookmarks.Add($bk, $rng) | Out-Null
    $rng.InsertAfter($text)
    $d.Bookmarks($bk).Delete()
    $rng.Collapse(0)
}

# ---------------------------------------------------------------------------
# 1. "I - DESCRIÇÃO" -> "I - " + "tipodocumento"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("I - DESCRIÇÃO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "I - "
$rng.Collapse(0)
Insert-AsNewRun $rng "tipodocumento"

# ---------------------------------------------------------------------------
# 2. "Vem ao exame desta Comissão o " -> "Vem ao exame desta Comissão o" +
#    " " + "pedidoAprovacao" + " "
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Vem ao exame desta Comissão o ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "Vem ao exame desta Comissão o"
$rng.Collapse(0)
Insert-AsNewRun $rng " "
Insert-AsNewRun $rng "pedidoAprovacao"
Insert-AsNewRun $rng " "

# ---------------------------------------------------------------------------
# 3. ", sem alocação de Carga Horária." -> ", " + "cargaHoraria" + "."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(", sem alocação de Carga Horária.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = ", "
$rng.Collapse(0)
Insert-AsNewRun $rng "cargaHoraria"
Insert-AsNewRun $rng "."
